# Scheduled-runner data refresh: updates cached market-price / profit
# figures (columns H-N) for specific leve rows across the ALC, ARM, BSM,
# CRP, CUL, GSM, LTW and WVR sheets. A few rows also gain or lose their
# LeveProfitNQ/HQ (M/N) cell depending on whether NQ or HQ pricing data
# was available for that refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 1259.7273
$ws.Range("I98").Value = 908.1429000000001
$ws.Range("K98").Value = 908.1429000000001
$ws.Range("M98").Value = 589.8570999999999

$ws.Range("H122").Value = 1259.7273
$ws.Range("I122").Value = 908.1429000000001
$ws.Range("K122").Value = 2724.4287
$ws.Range("M122").Value = -274.4287000000004

$ws.Range("H132").Value = 1551.4667
$ws.Range("I132").Value = 1681.4166
$ws.Range("K132").Value = 5044.2498
$ws.Range("M132").Value = -2514.2498

$ws.Range("H137").Value = 3827.3914
$ws.Range("I137").Value = 1798.8182
$ws.Range("J137").Value = 5686.9165
$ws.Range("K137").Value = 5396.4546
$ws.Range("L137").Value = 17060.7495
$ws.Range("M137").Value = -2846.4546
$ws.Range("N137").Value = -22160.7495

$ws.Range("H138").Value = 3151.8774
$ws.Range("I138").Value = 1353.5238
$ws.Range("J138").Value = 4500.643
$ws.Range("K138").Value = 4060.5714
$ws.Range("L138").Value = 13501.929
$ws.Range("M138").Value = 1079.4286
$ws.Range("N138").Value = -23781.929

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2323.8333
$ws.Range("I45").Value = 2088.6
$ws.Range("K45").Value = 2088.6
$ws.Range("M45").Value = -1711.6

$ws.Range("H61").Value = 2926.5
$ws.Range("I61").Value = 2912.4
$ws.Range("K61").Value = 2912.4
$ws.Range("M61").Value = -2700.4

$ws.Range("H122").Value = 2569.5715
$ws.Range("J122").Value = 1499
$ws.Range("L122").Value = 4497
$ws.Range("N122").Value = -9397

$ws.Range("H132").Value = 2545.64
$ws.Range("I132").Value = 2342.818
$ws.Range("K132").Value = 7028.454000000001
$ws.Range("M132").Value = -4498.454000000001

$ws.Range("H136").Value = 2926.5
$ws.Range("I136").Value = 2912.4
$ws.Range("K136").Value = 8737.200000000001
$ws.Range("M136").Value = -6187.200000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3799.125
$ws.Range("I20").Value = 4071.9092
$ws.Range("J20").Value = 3199
$ws.Range("K20").Value = 4071.9092
$ws.Range("L20").Value = 3199
$ws.Range("M20").Value = -3824.9092
$ws.Range("N20").Value = -3693

$ws.Range("H80").Value = 1031.4
$ws.Range("J80").Value = 1043.4
$ws.Range("L80").Value = 1043.4
$ws.Range("N80").Value = -3039.4

$ws.Range("H83").Value = 1031.4
$ws.Range("J83").Value = 1043.4
$ws.Range("L83").Value = 5217
$ws.Range("N83").Value = -15201

$ws.Range("H134").Value = 3861.7273
$ws.Range("I134").Value = 3861.7273
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 11585.1819
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -9050.1819
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1000
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 1000
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 1000
$ws.Range("N22").Value = -1700
$ws.Range("M22").ClearContents()

$ws.Range("H31").Value = 2247.3
$ws.Range("I31").Value = 2435
$ws.Range("K31").Value = 2435
$ws.Range("M31").Value = -2140

$ws.Range("H34").Value = 2247.3
$ws.Range("I34").Value = 2435
$ws.Range("K34").Value = 2435
$ws.Range("M34").Value = -2233

$ws.Range("H122").Value = 2136.3635
$ws.Range("I122").Value = 2200
$ws.Range("K122").Value = 6600
$ws.Range("M122").Value = -4150

$ws.Range("H134").Value = 1462
$ws.Range("I134").Value = 1541
$ws.Range("J134").Value = 948.5
$ws.Range("K134").Value = 4623
$ws.Range("L134").Value = 2845.5
$ws.Range("M134").Value = -2088
$ws.Range("N134").Value = -7915.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 796.6667
$ws.Range("I5").Value = 890
$ws.Range("J5").Value = 750
$ws.Range("K5").Value = 2670
$ws.Range("L5").Value = 2250
$ws.Range("M5").Value = -2558
$ws.Range("N5").Value = -2474

$ws.Range("H23").Value = 314.66666
$ws.Range("J23").Value = 322
$ws.Range("L23").Value = 966
$ws.Range("N23").Value = -1436

$ws.Range("H33").Value = 961
$ws.Range("I33").Value = 154.2
$ws.Range("K33").Value = 925.1999999999999
$ws.Range("M33").Value = -642.1999999999999

$ws.Range("H40").Value = 128.57143
$ws.Range("I40").Value = 128.57143
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 514.28572
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -445.28572
$ws.Range("N40").ClearContents()

$ws.Range("H56").Value = 19315.42
$ws.Range("I56").Value = 19315.42
$ws.Range("K56").Value = 19315.42
$ws.Range("M56").Value = -18785.42

$ws.Range("H68").Value = 1545.7142
$ws.Range("I68").Value = 1723.3334
$ws.Range("J68").Value = 1412.5
$ws.Range("K68").Value = 5170.0002
$ws.Range("L68").Value = 4237.5
$ws.Range("M68").Value = -4359.0002
$ws.Range("N68").Value = -5859.5

$ws.Range("H71").Value = 1545.7142
$ws.Range("I71").Value = 1723.3334
$ws.Range("J71").Value = 1412.5
$ws.Range("K71").Value = 15510.0006
$ws.Range("L71").Value = 12712.5
$ws.Range("M71").Value = -11454.0006
$ws.Range("N71").Value = -20824.5

$ws.Range("H109").Value = 2570.2856
$ws.Range("I109").Value = 1000
$ws.Range("J109").Value = 3198.4
$ws.Range("K109").Value = 3000
$ws.Range("L109").Value = 9595.200000000001
$ws.Range("M109").Value = -1960
$ws.Range("N109").Value = -11675.2

$ws.Range("H135").Value = 796.6667
$ws.Range("I135").Value = 890
$ws.Range("J135").Value = 750
$ws.Range("K135").Value = 8010
$ws.Range("L135").Value = 6750
$ws.Range("M135").Value = -5475
$ws.Range("N135").Value = -11820

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 463
$ws.Range("J107").Value = 493
$ws.Range("L107").Value = 493
$ws.Range("N107").Value = -4333

$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1074.75
$ws.Range("I22").Value = 999.5
$ws.Range("J22").Value = 1150
$ws.Range("K22").Value = 999.5
$ws.Range("L22").Value = 1150
$ws.Range("M22").Value = -704.5
$ws.Range("N22").Value = -1740

$ws.Range("H27").Value = 1074.75
$ws.Range("I27").Value = 999.5
$ws.Range("J27").Value = 1150
$ws.Range("K27").Value = 999.5
$ws.Range("L27").Value = 1150
$ws.Range("M27").Value = -892.5
$ws.Range("N27").Value = -1364

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 798.3333
$ws.Range("I107").Value = 798.5
$ws.Range("K107").Value = 2395.5
$ws.Range("M107").Value = -475.5

$ws.Range("H136").Value = 1481.2106
$ws.Range("J136").Value = 2742
$ws.Range("L136").Value = 8226
$ws.Range("N136").Value = -13326
